$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-09-09 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-10 Sunday", 2) | Out-Null
$d.Content.Find.Execute("72×70=5040", $true, $false, $false, $false, $false, $true, 1, $false, "56×27=1512", 2) | Out-Null
$d.Content.Find.Execute("85×44=3740", $true, $false, $false, $false, $false, $true, 1, $false, "55×59=3245", 2) | Out-Null
$d.Content.Find.Execute("49×21=1029", $true, $false, $false, $false, $false, $true, 1, $false, "60×18=1080", 2) | Out-Null
$d.Content.Find.Execute("49×30=1470", $true, $false, $false, $false, $false, $true, 1, $false, "63×44=2772", 2) | Out-Null
$d.Content.Find.Execute("32×69=2208", $true, $false, $false, $false, $false, $true, 1, $false, "37×61=2257", 2) | Out-Null
$d.Content.Find.Execute("41×34=1394", $true, $false, $false, $false, $false, $true, 1, $false, "80×13=1040", 2) | Out-Null
$d.Content.Find.Execute("53×13=689", $true, $false, $false, $false, $false, $true, 1, $false, "36×99=3564", 2) | Out-Null
$d.Content.Find.Execute("49×37=1813", $true, $false, $false, $false, $false, $true, 1, $false, "13×20=260", 2) | Out-Null
$d.Content.Find.Execute("12×43=516", $true, $false, $false, $false, $false, $true, 1, $false, "42×60=2520", 2) | Out-Null
$d.Content.Find.Execute("94×42=3948", $true, $false, $false, $false, $false, $true, 1, $false, "23×27=621", 2) | Out-Null
$d.Content.Find.Execute("64×54=3456", $true, $false, $false, $false, $false, $true, 1, $false, "23×22=506", 2) | Out-Null
$d.Content.Find.Execute("84×30=2520", $true, $false, $false, $false, $false, $true, 1, $false, "22×59=1298", 2) | Out-Null
$d.Content.Find.Execute("35×54=1890", $true, $false, $false, $false, $false, $true, 1, $false, "55×33=1815", 2) | Out-Null
$d.Content.Find.Execute("42×41=1722", $true, $false, $false, $false, $false, $true, 1, $false, "74×53=3922", 2) | Out-Null
$d.Content.Find.Execute("97×49=4753", $true, $false, $false, $false, $false, $true, 1, $false, "30×99=2970", 2) | Out-Null
$d.Content.Find.Execute("47×55=2585", $true, $false, $false, $false, $false, $true, 1, $false, "77×58=4466", 2) | Out-Null
$d.Content.Find.Execute("40×72=2880", $true, $false, $false, $false, $false, $true, 1, $false, "83×16=1328", 2) | Out-Null
$d.Content.Find.Execute("76×51=3876", $true, $false, $false, $false, $false, $true, 1, $false, "17×90=1530", 2) | Out-Null
$d.Content.Find.Execute("39×85=3315", $true, $false, $false, $false, $false, $true, 1, $false, "27×42=1134", 2) | Out-Null
$d.Content.Find.Execute("98×80=7840", $true, $false, $false, $false, $false, $true, 1, $false, "25×53=1325", 2) | Out-Null
$d.Content.Find.Execute("49×31=1519", $true, $false, $false, $false, $false, $true, 1, $false, "45×55=2475", 2) | Out-Null
$d.Content.Find.Execute("17×11=187", $true, $false, $false, $false, $false, $true, 1, $false, "76×69=5244", 2) | Out-Null
$d.Content.Find.Execute("43×70=3010", $true, $false, $false, $false, $false, $true, 1, $false, "66×41=2706", 2) | Out-Null
$d.Content.Find.Execute("59×33=1947", $true, $false, $false, $false, $false, $true, 1, $false, "43×68=2924", 2) | Out-Null
$d.Content.Find.Execute("41×40=1640", $true, $false, $false, $false, $false, $true, 1, $false, "44×88=3872", 2) | Out-Null
